# Auto-generated edit script: applies cell-level numeric updates
# to the Brynhildr_Profits workbook per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2546.5
$ws.Range("I18").Value = 2546.5
$ws.Range("K18").Value = 2546.5
$ws.Range("M18").Value = -2262.5

$ws.Range("H62").Value = 2016.8422
$ws.Range("I62").Value = 1645
$ws.Range("K62").Value = 1645
$ws.Range("M62").Value = -1021

$ws.Range("H65").Value = 2016.8422
$ws.Range("I65").Value = 1645
$ws.Range("K65").Value = 8225
$ws.Range("M65").Value = -5105

$ws.Range("H80").Value = 1682.0294
$ws.Range("I80").Value = 1455.3158
$ws.Range("K80").Value = 4365.9474
$ws.Range("M80").Value = -3367.9474

$ws.Range("H83").Value = 1682.0294
$ws.Range("I83").Value = 1455.3158
$ws.Range("K83").Value = 13097.8422
$ws.Range("M83").Value = -8105.842200000001

$ws.Range("H100").Value = 1675.1111
$ws.Range("I100").Value = 1494.9333
$ws.Range("K100").Value = 1494.9333
$ws.Range("M100").Value = -953.9332999999999

$ws.Range("H132").Value = 6216
$ws.Range("I132").Value = 6529.6
$ws.Range("K132").Value = 19588.8
$ws.Range("M132").Value = -17058.8

$ws.Range("H141").Value = 5101.0605
$ws.Range("I141").Value = 3060.0417
$ws.Range("K141").Value = 9180.125100000001
$ws.Range("M141").Value = -4000.125100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2087.0625
$ws.Range("I2").Value = 1888.909
$ws.Range("K2").Value = 1888.909
$ws.Range("M2").Value = -1775.909

$ws.Range("H116").Value = 2087.0625
$ws.Range("I116").Value = 1888.909
$ws.Range("K116").Value = 1888.909
$ws.Range("M116").Value = 405.0909999999999

$ws.Range("H121").Value = 44127.5
$ws.Range("J121").Value = 44127.5
$ws.Range("L121").Value = 44127.5
$ws.Range("N121").Value = -47621.5

$ws.Range("H122").Value = 1541.3334
$ws.Range("I122").Value = 1609
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 4827
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -2377
$ws.Range("N122").Value = -7900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2087.0625
$ws.Range("I3").Value = 1888.909
$ws.Range("K3").Value = 1888.909
$ws.Range("M3").Value = -1774.909

$ws.Range("H20").Value = 42781.8
$ws.Range("I20").Value = 65396.062
$ws.Range("J20").Value = 2578.6667
$ws.Range("K20").Value = 65396.062
$ws.Range("L20").Value = 2578.6667
$ws.Range("M20").Value = -65149.062
$ws.Range("N20").Value = -3072.6667

$ws.Range("H53").Value = 65390
$ws.Range("J53").Value = 65390
$ws.Range("L53").Value = 65390
$ws.Range("N53").Value = -66538

$ws.Range("H94").Value = 8217.105
$ws.Range("J94").Value = 8204.333000000001
$ws.Range("L94").Value = 8204.333000000001
$ws.Range("N94").Value = -9106.333000000001

$ws.Range("H99").Value = 7164.9443
$ws.Range("I99").Value = 9229.462
$ws.Range("K99").Value = 9229.462
$ws.Range("M99").Value = -7731.462

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws.Range("H134").Value = 800
$ws.Range("I134").Value = 800
$ws.Range("K134").Value = 2400
$ws.Range("M134").Value = 135

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 112125
$ws.Range("J28").Value = 49500
$ws.Range("L28").Value = 49500
$ws.Range("N28").Value = -49990

$ws.Range("H31").Value = 2605.6667
$ws.Range("I31").Value = 3496.75
$ws.Range("J31").Value = 1587.2858
$ws.Range("K31").Value = 3496.75
$ws.Range("L31").Value = 1587.2858
$ws.Range("M31").Value = -3201.75
$ws.Range("N31").Value = -2177.2858

$ws.Range("H34").Value = 2605.6667
$ws.Range("I34").Value = 3496.75
$ws.Range("J34").Value = 1587.2858
$ws.Range("K34").Value = 3496.75
$ws.Range("L34").Value = 1587.2858
$ws.Range("M34").Value = -3294.75
$ws.Range("N34").Value = -1991.2858

$ws.Range("H58").Value = 4901.5835
$ws.Range("I58").Value = 3282.3157
$ws.Range("K58").Value = 3282.3157
$ws.Range("M58").Value = -3079.3157

$ws.Range("H68").Value = 37566.668
$ws.Range("J68").Value = 37566.668
$ws.Range("L68").Value = 37566.668
$ws.Range("N68").Value = -39064.668

$ws.Range("H71").Value = 37566.668
$ws.Range("J71").Value = 37566.668
$ws.Range("L71").Value = 112700.004
$ws.Range("N71").Value = -120188.004

$ws.Range("H99").Value = 16563.2
$ws.Range("J99").Value = 2935.3333
$ws.Range("L99").Value = 2935.3333
$ws.Range("N99").Value = -5931.3333

$ws.Range("H105").Value = 21690
$ws.Range("I105").Value = 21690
$ws.Range("K105").Value = 21690
$ws.Range("M105").Value = -19943

$ws.Range("H124").Value = 67999
$ws.Range("J124").Value = 67999
$ws.Range("L124").Value = 67999
$ws.Range("N124").Value = -72909

$ws.Range("H126").Value = 16563.2
$ws.Range("J126").Value = 2935.3333
$ws.Range("L126").Value = 8805.999899999999
$ws.Range("N126").Value = -13745.9999

$ws.Range("H136").Value = 4901.5835
$ws.Range("I136").Value = 3282.3157
$ws.Range("K136").Value = 9846.947100000001
$ws.Range("M136").Value = -7296.947100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 639.5
$ws.Range("I2").Value = 99
$ws.Range("J2").Value = 747.6
$ws.Range("K2").Value = 594
$ws.Range("L2").Value = 4485.6
$ws.Range("M2").Value = -481
$ws.Range("N2").Value = -4711.6

$ws.Range("H124").Value = 13453.375
$ws.Range("I124").Value = 10875.667
$ws.Range("K124").Value = 32627.001
$ws.Range("M124").Value = -27717.001

$ws.Range("H129").Value = 2015.2858
$ws.Range("I129").Value = 1175.6666
$ws.Range("J129").Value = 2645
$ws.Range("K129").Value = 3526.9998
$ws.Range("L129").Value = 7935
$ws.Range("M129").Value = 1473.0002
$ws.Range("N129").Value = -17935

$ws.Range("H130").Value = 12037.571
$ws.Range("I130").Value = 4030
$ws.Range("K130").Value = 12090
$ws.Range("M130").Value = -7070

$ws.Range("H138").Value = 24865.74
$ws.Range("I138").Value = 40230.555
$ws.Range("J138").Value = 17183.334
$ws.Range("K138").Value = 120691.665
$ws.Range("L138").Value = 51550.00199999999
$ws.Range("M138").Value = -115551.665
$ws.Range("N138").Value = -61830.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 3960
$ws.Range("I23").Value = 2606
$ws.Range("J23").Value = 4862.6665
$ws.Range("K23").Value = 2606
$ws.Range("L23").Value = 4862.6665
$ws.Range("M23").Value = -2383
$ws.Range("N23").Value = -5308.6665

$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H122").Value = 3771.9167
$ws.Range("I122").Value = 4023.9092
$ws.Range("K122").Value = 12071.7276
$ws.Range("M122").Value = -9621.7276

$ws.Range("H126").Value = 2587.8096
$ws.Range("I126").Value = 2275.4546
$ws.Range("J126").Value = 2931.4
$ws.Range("K126").Value = 6826.3638
$ws.Range("L126").Value = 8794.200000000001
$ws.Range("M126").Value = -4356.3638
$ws.Range("N126").Value = -13734.2

$ws.Range("H132").Value = 9284.235000000001
$ws.Range("I132").Value = 9948.799999999999
$ws.Range("J132").Value = 4300
$ws.Range("K132").Value = 29846.4
$ws.Range("L132").Value = 12900
$ws.Range("M132").Value = -27316.4
$ws.Range("N132").Value = -17960

$ws.Range("H136").Value = 57996
$ws.Range("J136").Value = 57996
$ws.Range("L136").Value = 173988
$ws.Range("N136").Value = -179088

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3881.3684
$ws.Range("I40").Value = 3487.625
$ws.Range("K40").Value = 3487.625
$ws.Range("M40").Value = -3351.625

$ws.Range("H61").Value = 9093.526
$ws.Range("I61").Value = 8463.134
$ws.Range("K61").Value = 8463.134
$ws.Range("M61").Value = -8261.134

$ws.Range("H113").Value = 9093.526
$ws.Range("I113").Value = 8463.134
$ws.Range("K113").Value = 8463.134
$ws.Range("M113").Value = -6293.134

$ws.Range("H122").Value = 6038.625
$ws.Range("I122").Value = 5051.5
$ws.Range("K122").Value = 15154.5
$ws.Range("M122").Value = -12704.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 60000
$ws.Range("J5").Value = 60000
$ws.Range("L5").Value = 60000
$ws.Range("N5").Value = -60224

$ws.Range("H107").Value = 1732.8684
$ws.Range("I107").Value = 1043.2307
$ws.Range("K107").Value = 3129.6921
$ws.Range("M107").Value = -1209.6921

